$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J1: was shared string "r" -> now numeric 0.6
$ws.Range("J1").Value = 0.6

# K1: was shared string "s" -> now numeric 0.5
$ws.Range("K1").Value = 0.5

# K2:K51: was 0.6 -> now 0.5
$ws.Range("K2:K51").Value = 0.5

# Update sheet view: topLeftCell A39 -> A40, selection J2:K51 -> K52:K54
$excel.ActiveWindow.ScrollRow = 40
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("K52:K54").Select()
